# Apply the "save pw as string" schedule update to the Full Schedule sheet.
# The block of cells in columns E:G (the left mini-roster) and the block in
# columns I:K (the right mini-roster) both had their leading blank/placeholder
# rows removed, which shows up as the surviving rows 98-116 inheriting values
# from further down (E:G shifted up 3 rows, I:J:K shifted up 1 row), and the
# three now-unused trailing rows (117-119) being deleted outright.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F / J notes that slide up into rows 98-103 ---
$ws.Range("F98").Value = "IL: 4:15 AM MEET AT OFFICE"
$ws.Range("J98").Value = "IL: 4:00 AM MEET AT OFFICE" + [char]10 + "MD: 4:45 MEET HOME DEPOT VERONA RD"

$ws.Range("F99").Value = "6:00 AM START"

$ws.Range("F100").Value = "DC5-FINANCIAL"

$ws.Range("F104").Value = ""

# --- Left roster (E/F/G) entries shift up, renumbering 1)-12) ---
$ws.Range("E105").Value = "1)"
$ws.Range("G105").Value = "@ Store, Equip"

$ws.Range("E106").Value = "2)"
$ws.Range("G106").Value = "Driver, 1/2" + [char]10 + "Silver Van"

$ws.Range("E107").Value = "3)"
$ws.Range("F107").Value = "Anisha"

$ws.Range("G109").Value = ""
$ws.Range("J109").Value = ""

# --- Right roster (I/J/K) entries shift up, renumbering 1)-6) ---
$ws.Range("I110").Value = "1)"
$ws.Range("J110").Value = "Sarah"
$ws.Range("K110").Value = "Driver," + [char]10 + "Equip"

$ws.Range("G111").Value = ""
$ws.Range("I111").Value = "2)"
$ws.Range("J111").Value = "Ashley P"
$ws.Range("K111").Value = ""

$ws.Range("K112").Value = ""

$ws.Range("G114").Value = "Driver, 1/2" + [char]10 + "Gray Van"

$ws.Range("G115").Value = "Driver, 1/2" + [char]10 + "Gray Van"
$ws.Range("K115").Value = "Driver, 1/2" + [char]10 + "Silver Van"

$ws.Range("I116").Value = ""
$ws.Range("J116").Value = ""
$ws.Range("K116").Value = ""

# --- The trailing rows that no longer hold any data are removed outright ---
$ws.Rows("117:119").Delete()
